# Update NATMI TPM output (LR-pairs_lrc2p F2-Gp9) for rows 2-21, columns D (Target
# cluster) through T (Edge total expression derived specificity) with refreshed
# values from the new TPM-based computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",3.0,1.0,0.8966636666666666,2.689991,0.334725143386341,0.3647360854412732,1.0,0.3333333333333333,0.2231716666666667,0.669515,0.02293167831634903,0.02293167831634903,0.2001099249294444,1.800989324365,0.007675809312529375,0.008364010581703671),
    @("FAPs",3.0,1.0,0.8966636666666666,2.689991,0.334725143386341,0.3647360854412732,1.0,0.3333333333333333,0.223816,0.671448,0.02299788584595704,0.02299788584595704,0.2006876752186667,1.806189076968,0.007697970637370671,0.008388158856879633),
    @("Inflammatory-Mac",3.0,1.0,0.8966636666666666,2.689991,0.334725143386341,0.3647360854412732,3.0,1.0,4.462895666666667,13.388687,0.4585783191747521,0.4585783191747521,4.001716392424112,36.015447531817,0.1534976936396362,0.1672600610040378),
    @("Resolving-Mac",3.0,1.0,0.8966636666666666,2.689991,0.334725143386341,0.3647360854412732,3.0,1.0,4.822141666666667,14.466425,0.4954921166629419,0.4954921166629419,4.323839228019445,38.914553052175,0.1658536697968048,0.1807238549986521),
    @("ECs",3.0,1.0,0.6827986666666667,2.048396,0.2548891965854188,0.2777421703171357,1.0,0.3333333333333333,0.2231716666666667,0.669515,0.02293167831634903,0.02293167831634903,0.1523813164377778,1.37143184794,0.005845037062409473,0.006369094104597179),
    @("FAPs",3.0,1.0,0.6827986666666667,2.048396,0.2548891965854188,0.2777421703171357,1.0,0.3333333333333333,0.223816,0.671448,0.02299788584595704,0.02299788584595704,0.1528212663786667,1.375391397408,0.005861912646439164,0.006387482727561844),
    @("Inflammatory-Mac",3.0,1.0,0.6827986666666667,2.048396,0.2548891965854188,0.2777421703171357,3.0,1.0,4.462895666666667,13.388687,0.4585783191747521,0.4585783191747521,3.047259210672445,27.425332896052,0.1168866593459443,0.1273665376279798),
    @("Resolving-Mac",3.0,1.0,0.6827986666666667,2.048396,0.2548891965854188,0.2777421703171357,3.0,1.0,4.822141666666667,14.466425,0.4954921166629419,0.4954921166629419,3.292551900477778,29.6329671043,0.1262955875306259,0.1376190558569969),
    @("ECs",2.0,0.6666666666666666,0.221369,0.664107,0.08263719499391366,0.09004631892602898,1.0,0.3333333333333333,0.2231716666666667,0.669515,0.02293167831634903,0.02293167831634903,0.04940328867833333,0.444629598105,0.001895009572565836,0.002064913219183068),
    @("FAPs",2.0,0.6666666666666666,0.221369,0.664107,0.08263719499391366,0.09004631892602898,1.0,0.3333333333333333,0.223816,0.671448,0.02299788584595704,0.02299788584595704,0.04954592410400001,0.4459133169360001,0.001900480777100119,0.002070874963509455),
    @("Inflammatory-Mac",2.0,0.6666666666666666,0.221369,0.664107,0.08263719499391366,0.09004631892602898,3.0,1.0,4.462895666666667,13.388687,0.4585783191747521,0.4585783191747521,0.9879467508343335,8.891520757509001,0.03789562598162517,0.04129328958097204),
    @("Resolving-Mac",2.0,0.6666666666666666,0.221369,0.664107,0.08263719499391366,0.09004631892602898,3.0,1.0,4.822141666666667,14.466425,0.4954921166629419,0.4954921166629419,1.067472678608334,9.607254107475,0.04094607866262254,0.04461724116236442),
    @("ECs",2.0,1.0,0.6612465000000001,1.322493,0.2468437584284291,0.1793169270244717,1.0,0.3333333333333333,0.2231716666666667,0.669515,0.02293167831634903,0.02293167831634903,0.1475714834825,0.885428900895,0.005660541662679305,0.00411203808720142),
    @("FAPs",2.0,1.0,0.6612465000000001,1.322493,0.2468437584284291,0.1793169270244717,1.0,0.3333333333333333,0.223816,0.671448,0.02299788584595704,0.02299788584595704,0.147997546644,0.8879852798640001,0.005676884578124008,0.004123910217956609),
    @("Inflammatory-Mac",2.0,1.0,0.6612465000000001,1.322493,0.2468437584284291,0.1793169270244717,3.0,1.0,4.462895666666667,13.388687,0.4585783191747521,0.4585783191747521,2.9510741394485,17.706444836691,0.1131971958388876,0.08223085499446393),
    @("Resolving-Mac",2.0,1.0,0.6612465000000001,1.322493,0.2468437584284291,0.1793169270244717,3.0,1.0,4.822141666666667,14.466425,0.4954921166629419,0.4954921166629419,3.1886242995875,19.131745797525,0.1223091363487382,0.08885012372484978),
    @("ECs",3.0,1.0,0.216728,0.650184,0.08090470660589748,0.08815849829109049,1.0,0.3333333333333333,0.2231716666666667,0.669515,0.02293167831634903,0.02293167831634903,0.04836754897333333,0.43530794076,0.001855280706165039,0.002021622323663692),
    @("FAPs",3.0,1.0,0.216728,0.650184,0.08090470660589748,0.08815849829109049,1.0,0.3333333333333333,0.223816,0.671448,0.02299788584595704,0.02299788584595704,0.048507194048,0.436564746432,0.001860637206923076,0.002027459080049498),
    @("Inflammatory-Mac",3.0,1.0,0.216728,0.650184,0.08090470660589748,0.08815849829109049,3.0,1.0,4.462895666666667,13.388687,0.4585783191747521,0.4585783191747521,0.9672344520453334,8.705110068408,0.03710114436865893,0.04042757596729853),
    @("Resolving-Mac",3.0,1.0,0.216728,0.650184,0.08090470660589748,0.08815849829109049,3.0,1.0,4.822141666666667,14.466425,0.4954921166629419,0.4954921166629419,1.045093119133333,9.4058380722,0.04008764432415043,0.04368184092007877)
)

$startRow = 2
$startCol = 4  # column D

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $startCol + $j
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
